# Adds "Sheet2" (cron-generated rsi / macd / atr order-decision matrix) after Sheet1,
# matching the commit "Added rsi, macd and atr in cron".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet goes right after Sheet1 and becomes the active tab.
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "Sheet2"

# --- Priming assignments: establish shared-string insertion order ---
$ws.Range("A1").Value = "momentum_rsi"
$ws.Range("B1").Value = "trend_macd"
$ws.Range("C1").Value = "volatility_atr"
$ws.Range("A2").Value = "buy"
$ws.Range("C11").Value = "low"
$ws.Range("C2").Value = "high"
$ws.Range("A3").Value = "sell"
$ws.Range("D1").Value = "no_order"
$ws.Range("E1").Value = "buy_order"
$ws.Range("F1").Value = "sell_order"
$ws.Range("E2").Value = "hold"
$ws.Range("D3").Value = "none"
$ws.Range("F2").Value = "exit"

# --- Remaining grid values (reuse already-interned shared strings) ---
$ws.Cells.Item(2, 2).Value = "buy"
$ws.Cells.Item(2, 4).Value = "buy"
$ws.Cells.Item(3, 2).Value = "buy"
$ws.Cells.Item(3, 3).Value = "high"
$ws.Cells.Item(3, 5).Value = "exit"
$ws.Cells.Item(3, 6).Value = "exit"
$ws.Cells.Item(4, 1).Value = "neutral"
$ws.Cells.Item(4, 2).Value = "buy"
$ws.Cells.Item(4, 3).Value = "high"
$ws.Cells.Item(4, 4).Value = "none"
$ws.Cells.Item(4, 5).Value = "hold"
$ws.Cells.Item(4, 6).Value = "exit"
$ws.Cells.Item(5, 1).Value = "buy"
$ws.Cells.Item(5, 2).Value = "sell"
$ws.Cells.Item(5, 3).Value = "high"
$ws.Cells.Item(5, 4).Value = "none"
$ws.Cells.Item(5, 5).Value = "exit"
$ws.Cells.Item(5, 6).Value = "exit"
$ws.Cells.Item(6, 1).Value = "sell"
$ws.Cells.Item(6, 2).Value = "sell"
$ws.Cells.Item(6, 3).Value = "high"
$ws.Cells.Item(6, 4).Value = "buy"
$ws.Cells.Item(6, 5).Value = "exit"
$ws.Cells.Item(6, 6).Value = "hold"
$ws.Cells.Item(7, 1).Value = "neutral"
$ws.Cells.Item(7, 2).Value = "sell"
$ws.Cells.Item(7, 3).Value = "high"
$ws.Cells.Item(7, 4).Value = "none"
$ws.Cells.Item(7, 5).Value = "exit"
$ws.Cells.Item(7, 6).Value = "hold"
$ws.Cells.Item(8, 1).Value = "buy"
$ws.Cells.Item(8, 2).Value = "neutral"
$ws.Cells.Item(8, 3).Value = "high"
$ws.Cells.Item(8, 4).Value = "none"
$ws.Cells.Item(8, 5).Value = "exit"
$ws.Cells.Item(8, 6).Value = "exit"
$ws.Cells.Item(9, 1).Value = "sell"
$ws.Cells.Item(9, 2).Value = "neutral"
$ws.Cells.Item(9, 3).Value = "high"
$ws.Cells.Item(9, 4).Value = "none"
$ws.Cells.Item(9, 5).Value = "exit"
$ws.Cells.Item(9, 6).Value = "exit"
$ws.Cells.Item(10, 1).Value = "neutral"
$ws.Cells.Item(10, 2).Value = "neutral"
$ws.Cells.Item(10, 3).Value = "high"
$ws.Cells.Item(10, 4).Value = "none"
$ws.Cells.Item(10, 5).Value = "exit"
$ws.Cells.Item(10, 6).Value = "exit"
$ws.Cells.Item(11, 1).Value = "buy"
$ws.Cells.Item(11, 2).Value = "buy"
$ws.Cells.Item(11, 4).Value = "none"
$ws.Cells.Item(11, 5).Value = "exit"
$ws.Cells.Item(11, 6).Value = "exit"
$ws.Cells.Item(12, 1).Value = "sell"
$ws.Cells.Item(12, 2).Value = "buy"
$ws.Cells.Item(12, 3).Value = "low"
$ws.Cells.Item(12, 4).Value = "none"
$ws.Cells.Item(12, 5).Value = "exit"
$ws.Cells.Item(12, 6).Value = "exit"
$ws.Cells.Item(13, 1).Value = "neutral"
$ws.Cells.Item(13, 2).Value = "buy"
$ws.Cells.Item(13, 3).Value = "low"
$ws.Cells.Item(13, 4).Value = "none"
$ws.Cells.Item(13, 5).Value = "exit"
$ws.Cells.Item(13, 6).Value = "exit"
$ws.Cells.Item(14, 1).Value = "buy"
$ws.Cells.Item(14, 2).Value = "sell"
$ws.Cells.Item(14, 3).Value = "low"
$ws.Cells.Item(14, 4).Value = "none"
$ws.Cells.Item(14, 5).Value = "exit"
$ws.Cells.Item(14, 6).Value = "exit"
$ws.Cells.Item(15, 1).Value = "sell"
$ws.Cells.Item(15, 2).Value = "sell"
$ws.Cells.Item(15, 3).Value = "low"
$ws.Cells.Item(15, 4).Value = "none"
$ws.Cells.Item(15, 5).Value = "exit"
$ws.Cells.Item(15, 6).Value = "exit"
$ws.Cells.Item(16, 1).Value = "neutral"
$ws.Cells.Item(16, 2).Value = "sell"
$ws.Cells.Item(16, 3).Value = "low"
$ws.Cells.Item(16, 4).Value = "none"
$ws.Cells.Item(16, 5).Value = "exit"
$ws.Cells.Item(16, 6).Value = "exit"
$ws.Cells.Item(17, 1).Value = "buy"
$ws.Cells.Item(17, 2).Value = "neutral"
$ws.Cells.Item(17, 3).Value = "low"
$ws.Cells.Item(17, 4).Value = "none"
$ws.Cells.Item(17, 5).Value = "exit"
$ws.Cells.Item(17, 6).Value = "exit"
$ws.Cells.Item(18, 1).Value = "sell"
$ws.Cells.Item(18, 2).Value = "neutral"
$ws.Cells.Item(18, 3).Value = "low"
$ws.Cells.Item(18, 4).Value = "none"
$ws.Cells.Item(18, 5).Value = "exit"
$ws.Cells.Item(18, 6).Value = "exit"
$ws.Cells.Item(19, 1).Value = "neutral"
$ws.Cells.Item(19, 2).Value = "neutral"
$ws.Cells.Item(19, 3).Value = "low"
$ws.Cells.Item(19, 4).Value = "none"
$ws.Cells.Item(19, 5).Value = "exit"
$ws.Cells.Item(19, 6).Value = "exit"

# Column widths (character units) for A:D - closest achievable values given
# this COM host's 1/6-character rounding granularity.
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.5
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 9

# Selection / active cell on the new sheet.
$ws.Range("G19").Select() | Out-Null
